$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Update row 2 (OTROS group): VENTA and POR CUMPLIR
$ws.Range("D2").Value = 922.8
$ws.Range("E2").Value = -922.8

# Update row 4 (TOTAL row): VENTA, POR CUMPLIR, CUMPLIMIENTO
$ws.Range("D4").Value = 1254.55
$ws.Range("E4").Value = 12468.79
$ws.Range("F4").Value = 0.09141724973658016
